$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 3 (pushes existing rows 3.. down to 6..)
$ws.Rows("3:5").Insert()

# Write new-string cells first, in the order the new shared strings should be
# appended (matches target sharedStrings uniqueCount order):
#   20 = "Getting-Started"                                  (C4)
#   21 = "Removed from sidebar; duplicated sidebar links."  (D4)
#   22 = "Installing-FFmpeg"                                (C3)
$ws.Range("C4").Value = "Getting-Started"
$ws.Range("D4").Value = "Removed from sidebar; duplicated sidebar links."
$ws.Range("C3").Value = "Installing-FFmpeg"

# New row 3: adapt_authoring | Installing-FFmpeg | Added new.
$ws.Range("A3").Value = 42104
$ws.Range("B3").Value = "adapt_authoring"
$ws.Range("D3").Value = "Added new."

# New row 4: adapt_authoring | Getting-Started | Removed from sidebar; duplicated sidebar links.
$ws.Range("A4").Value = 42104
$ws.Range("B4").Value = "adapt_authoring"

# New row 5: adapt_authoring | _Sidebar | Added new.
$ws.Range("A5").Value = 42104
$ws.Range("B5").Value = "adapt_authoring"
$ws.Range("C5").Value = "_Sidebar"
$ws.Range("D5").Value = "Added new."

# Match formatting (number format, styles) of the data rows below
$ws.Range("A3:D5").NumberFormat = $ws.Range("A6").NumberFormat

$ws.Range("D8").Select()
